# ============================================================================
# Edit: "added short circuit angle example"
# - Adds a new worksheet "FaultsPOC" between RVC and LineCodes
# - Adds seven new defined names used by the new sheet
# - Updates several inputs/values on RVC (line lengths, DP) which ripple
#   through existing formulas
# - Adds two helper columns (N/O) on RVC summarizing POC/Sub values
# - Builds out the FaultsPOC sheet: short circuit current phasor addition
#   example (magnitude/phasor sums, error comparisons)
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. RVC sheet: update input values (line lengths, DP) -------------------
# ----------------------------------------------------------------------
$rvc = $wb.Worksheets.Item("RVC")

# Segment lengths (column D) feed E/F/J/K formulas already on the sheet.
$rvc.Range("D4").Value = 2000
$rvc.Range("D5").Value = 0
$rvc.Range("D6").Value = 0

# Rewrite J/K formulas (now one consistent per-row H/I * $D / 5280 pattern)
$rvc.Range("J4").Formula = '=H4*$D4/5280'
$rvc.Range("K4").Formula = '=I4*$D4/5280'
$rvc.Range("J5").Formula = '=H5*$D5/5280'
$rvc.Range("K5").Formula = '=I5*$D5/5280'
$rvc.Range("J6").Formula = '=H6*$D6/5280'
$rvc.Range("K6").Formula = '=I6*$D6/5280'

# DP input (B9) changes from 922000 to 5000000
$rvc.Range("B9").Value = 5000000

# ----------------------------------------------------------------------
# 2. RVC sheet: new helper columns N/O (POC / Sub labels + values) ------
# ----------------------------------------------------------------------
$rvc.Range("N3").Value = "POC"
$rvc.Range("N3").Font.Bold = $true
$rvc.Range("O3").Value = "Sub"

$rvc.Range("N4").Value = 6600
$rvc.Range("O4").Value = 2000
$rvc.Range("N5").Value = 3960
$rvc.Range("O5").Value = 0
$rvc.Range("N6").Value = 1320
$rvc.Range("O6").Value = 0

$rvc.Range("N9").Value = 922000
$rvc.Range("O9").Value = 5000000
$rvc.Range("N9").NumberFormat = "0.00E+00"
$rvc.Range("O9").NumberFormat = "0.00E+00"

# ----------------------------------------------------------------------
# 3. Insert the new "FaultsPOC" worksheet between RVC and LineCodes -----
# ----------------------------------------------------------------------
$lineCodes = $wb.Worksheets.Item("LineCodes")
$faults = $wb.Worksheets.Add($null, $rvc)
$faults.Name = "FaultsPOC"

# ----------------------------------------------------------------------
# 4. New defined names living on FaultsPOC -------------------------------
# ----------------------------------------------------------------------
$wb.Names.Add('VTf', '=FaultsPOC!$B$1')
$wb.Names.Add('Sder', '=FaultsPOC!$B$2')
$wb.Names.Add('IscMach', '=FaultsPOC!$H$1')
$wb.Names.Add('IscIBR', '=FaultsPOC!$H$2')
$wb.Names.Add('IbaseDER', '=FaultsPOC!$H$3')
$wb.Names.Add('MachAng', '=FaultsPOC!$I$1')
$wb.Names.Add('IBRAng', '=FaultsPOC!$I$2')

# ----------------------------------------------------------------------
# 5. FaultsPOC content ----------------------------------------------------
# ----------------------------------------------------------------------

# Row 1 -------------------------------------------------------------
$faults.Range("A1").Value = "VT [V]"
$faults.Range("B1").Value = 12470
$faults.Range("G1").Value = "Mach Isc"
$faults.Range("H1").Value = 5
$faults.Range("I1").Value = -85
$faults.Range("M1").Value = "At POC"
$faults.Range("Q1").Value = "At Sub"

# Row 2 -------------------------------------------------------------
$faults.Range("A2").Value = "S [VA]"
$faults.Range("B2").Value = 5000000
$faults.Range("D2").Value = "Mag"
$faults.Range("E2").Value = "Angle"
$faults.Range("G2").Value = "IBR Isc"
$faults.Range("H2").Value = 1.2
$faults.Range("I2").Value = 0
$faults.Range("M2").Value = 922000
$faults.Range("Q2").Value = 5000000

# Row 3 -------------------------------------------------------------
$faults.Range("A3").Value = "Total"
$faults.Range("B3").Value = 0.19290196515151514
$faults.Range("C3").Value = 0.86139753939393937
$faults.Range("D3").Formula = '=SQRT(B3*B3+C3*C3)'
$faults.Range("E3").Formula = '=DEGREES(ATAN2(B3,C3))'
$faults.Range("G3").Value = "Ibase DER"
$faults.Range("H3").Formula = '=Sder/SQRT(3)/VTf'
$faults.Range("M3").Value = 3.08342545
$faults.Range("N3").Value = 2.2573286000000001
$faults.Range("Q3").Value = 0.19290196515151514
$faults.Range("R3").Value = 0.86139753939393937

# Row 4 -------------------------------------------------------------
$faults.Range("A4").Value = "Increment"
$faults.Range("B4").Value = 0.31714438939393941
$faults.Range("C4").Value = 1.2826096606060606
$faults.Range("D4").Formula = '=SQRT(B4*B4+C4*C4)'
$faults.Range("E4").Formula = '=DEGREES(ATAN2(B4,C4))'
$faults.Range("M4").Value = 4.16042545
$faults.Range("N4").Value = 5.4843286000000004
$faults.Range("Q4").Value = 0.31714438939393941
$faults.Range("R4").Value = 1.2826096606060606

# Row 5 -------------------------------------------------------------
$faults.Range("A5").Value = "Zslgf [" + [char]0x03A9 + "]"
$faults.Range("B5").Formula = '=2*(B3+B4)/3'
$faults.Range("C5").Formula = '=2*(C3+C4)/3'
$faults.Range("D5").Formula = '=SQRT(B5*B5+C5*C5)'
$faults.Range("E5").Formula = '=DEGREES(ATAN2(B5,C5))'

# Row 6 -------------------------------------------------------------
$faults.Range("B6").Value = "Mag"
$faults.Range("C6").Value = "Angle"
$faults.Range("D6").Value = "Re"
$faults.Range("E6").Value = "Im"

# Row 7 -------------------------------------------------------------
$faults.Range("A7").Value = "Islgf [A]"
$faults.Range("B7").Formula = '=VTf/SQRT(3)/D5'
$faults.Range("C7").Formula = '=-E5'
$faults.Range("D7").Formula = '=$B7*COS(RADIANS($C7))'
$faults.Range("E7").Formula = '=$B7*SIN(RADIANS($C7))'

# Row 8 -------------------------------------------------------------
$faults.Range("A8").Value = "Imach [A]"
$faults.Range("B8").Formula = '=IscMach*IbaseDER'
$faults.Range("C8").Formula = '=MachAng'
$faults.Range("D8").Formula = '=$B8*COS(RADIANS($C8))'
$faults.Range("E8").Formula = '=$B8*SIN(RADIANS($C8))'

# Row 9 -------------------------------------------------------------
$faults.Range("A9").Value = "Iibr [A]"
$faults.Range("B9").Formula = '=IscIBR*IbaseDER'
$faults.Range("C9").Formula = '=IBRAng'
$faults.Range("D9").Formula = '=$B9*COS(RADIANS($C9))'
$faults.Range("E9").Formula = '=$B9*SIN(RADIANS($C9))'

# Row 11 ------------------------------------------------------------
$faults.Range("A11").Value = "Adding Current Magnitudes:"

# Row 12 ------------------------------------------------------------
$faults.Range("A12").Value = "Machine"
$faults.Range("B12").Formula = '=B7+B8'

# Row 13 ------------------------------------------------------------
$faults.Range("A13").Value = "IBR"
$faults.Range("B13").Formula = '=B7+B9'

# Row 15 ------------------------------------------------------------
$faults.Range("A15").Value = "Adding Phasors:"

# Row 16 ------------------------------------------------------------
$faults.Range("B16").Value = "Re"
$faults.Range("C16").Value = "Im"
$faults.Range("D16").Value = "Mag"
$faults.Range("E16").Value = "Angle"

# Row 17 ------------------------------------------------------------
$faults.Range("A17").Value = "Machine"
$faults.Range("B17").Formula = '=D7+D8'
$faults.Range("C17").Formula = '=E7+E8'
$faults.Range("D17").Formula = '=SQRT(B17*B17+C17*C17)'
$faults.Range("E17").Formula = '=DEGREES(ATAN2(B17,C17))'

# Row 18 ------------------------------------------------------------
$faults.Range("A18").Value = "IBR"
$faults.Range("B18").Formula = '=D7+D9'
$faults.Range("C18").Formula = '=E7+E9'
$faults.Range("D18").Formula = '=SQRT(B18*B18+C18*C18)'
$faults.Range("E18").Formula = '=DEGREES(ATAN2(B18,C18))'

# Row 20 ------------------------------------------------------------
$faults.Range("A20").Value = "Magnitude Errors:"

# Row 21 ------------------------------------------------------------
$faults.Range("B21").Value = "Total"
$faults.Range("C21").Value = "Increment"

# Row 22 ------------------------------------------------------------
$faults.Range("A22").Value = "Machine"
$faults.Range("B22").Formula = '=B12/D17-1'
$faults.Range("C22").Formula = '=(B12-B7)/(D17-B7)-1'

# Row 23 ------------------------------------------------------------
$faults.Range("A23").Value = "IBR"
$faults.Range("B23").Formula = '=B13/D18-1'
$faults.Range("C23").Formula = '=(B13-B7)/(D18-B7)-1'

# ----------------------------------------------------------------------
# 6. Styling -------------------------------------------------------------
# ----------------------------------------------------------------------

# Bold section headers / labels
$boldCells = @("A1", "A2", "A3", "A4", "A11", "A15", "A20", "M1", "Q1")
foreach ($addr in $boldCells) {
    $faults.Range($addr).Font.Bold = $true
}

# Right-aligned column headers
$rightAlignCells = @("D2", "E2", "B6", "C6", "D6", "E6", "D16", "E16", "B21", "C21")
foreach ($addr in $rightAlignCells) {
    $faults.Range($addr).HorizontalAlignment = -4152   # xlRight
}

# Highlighted (yellow) four-decimal inputs
$yellowCells = @("B3", "C3", "B4", "C4")
foreach ($addr in $yellowCells) {
    $faults.Range($addr).NumberFormat = "0.0000"
    $faults.Range($addr).Interior.Color = 65535
}

# Two-decimal angle value
$faults.Range("C8").NumberFormat = "0.00"

# Percent-with-fill error cells (orange fill, 2-decimal percent)
$errCells = @("B22", "C22", "B23", "C23")
foreach ($addr in $errCells) {
    $faults.Range($addr).NumberFormat = "0.00%"
    $faults.Range($addr).Interior.Color = 49407
}

# Scientific notation for the POC/Sub summary values
$faults.Range("M2").NumberFormat = "0.00E+00"
$faults.Range("Q2").NumberFormat = "0.00E+00"

$faults.Range("A12").Select() | Out-Null

Write-Host "Edit applied"
